$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Fuel"
$ws.Range("D14").Value = "A0"
$ws.Range("A15").Value = "Temp "
$ws.Range("D15").Value = "A1"
$ws.Range("D16").Value = "A2"
$ws.Range("D17").Value = "A3"
$ws.Range("A18").Value = "Oil Temp"
$ws.Range("D18").Value = "A4"

$ws.Range("F26").Select()
